$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1600.25
$ws.Range("I40").Value = 2001
$ws.Range("K40").Value = 2001
$ws.Range("M40").Value = -1826
$ws.Range("H53").Value = 840
$ws.Range("I53").Value = 750
$ws.Range("J53").Value = 862.5
$ws.Range("K53").Value = 750
$ws.Range("L53").Value = 862.5
$ws.Range("M53").Value = -113
$ws.Range("N53").Value = -2136.5
$ws.Range("H62").Value = 5686.4287
$ws.Range("I62").Value = 1951.25
$ws.Range("J62").Value = 10666.667
$ws.Range("K62").Value = 1951.25
$ws.Range("L62").Value = 10666.667
$ws.Range("M62").Value = -1327.25
$ws.Range("N62").Value = -11914.667
$ws.Range("H64").Value = 2708.8235
$ws.Range("I64").Value = 2700
$ws.Range("J64").Value = 2718.75
$ws.Range("K64").Value = 2700
$ws.Range("L64").Value = 2718.75
$ws.Range("M64").Value = -2452
$ws.Range("N64").Value = -3214.75
$ws.Range("H65").Value = 5686.4287
$ws.Range("I65").Value = 1951.25
$ws.Range("J65").Value = 10666.667
$ws.Range("K65").Value = 9756.25
$ws.Range("L65").Value = 53333.335
$ws.Range("M65").Value = -6636.25
$ws.Range("N65").Value = -59573.335
$ws.Range("H67").Value = 2708.8235
$ws.Range("I67").Value = 2700
$ws.Range("J67").Value = 2718.75
$ws.Range("K67").Value = 2700
$ws.Range("L67").Value = 2718.75
$ws.Range("M67").Value = -1842
$ws.Range("N67").Value = -4434.75
$ws.Range("H121").Value = 888.7727
$ws.Range("J121").Value = 888.7727
$ws.Range("L121").Value = 2666.3181
$ws.Range("N121").Value = -6160.3181
$ws.Range("H123").Value = 42780
$ws.Range("J123").Value = 42780
$ws.Range("L123").Value = 42780
$ws.Range("N123").Value = -52580
$ws.Range("H129").Value = 73948.53999999999
$ws.Range("I129").Value = 316.33334
$ws.Range("J129").Value = 96038.2
$ws.Range("K129").Value = 949.0000200000001
$ws.Range("L129").Value = 288114.6
$ws.Range("M129").Value = 4050.99998
$ws.Range("N129").Value = -298114.6
$ws.Range("H132").Value = 26424510
$ws.Range("I132").Value = 32389850
$ws.Range("J132").Value = 6570
$ws.Range("K132").Value = 97169550
$ws.Range("L132").Value = 19710
$ws.Range("M132").Value = -97167020
$ws.Range("N132").Value = -24770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5212.045
$ws.Range("I32").Value = 2995.6233
$ws.Range("K32").Value = 2995.6233
$ws.Range("M32").Value = -2708.6233
$ws.Range("H110").Value = 895.03705
$ws.Range("I110").Value = 869.6667
$ws.Range("J110").Value = 983.8333
$ws.Range("K110").Value = 869.6667
$ws.Range("L110").Value = 983.8333
$ws.Range("M110").Value = 1175.3333
$ws.Range("N110").Value = -5073.8333
$ws.Range("H132").Value = 2896.525
$ws.Range("I132").Value = 2117.9583
$ws.Range("J132").Value = 4064.375
$ws.Range("K132").Value = 6353.874899999999
$ws.Range("L132").Value = 12193.125
$ws.Range("M132").Value = -3823.874899999999
$ws.Range("N132").Value = -17253.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 226.7037
$ws.Range("I80").Value = 107.71429
$ws.Range("J80").Value = 268.35
$ws.Range("K80").Value = 107.71429
$ws.Range("L80").Value = 268.35
$ws.Range("M80").Value = 890.28571
$ws.Range("N80").Value = -2264.35
$ws.Range("H83").Value = 226.7037
$ws.Range("I83").Value = 107.71429
$ws.Range("J83").Value = 268.35
$ws.Range("K83").Value = 538.57145
$ws.Range("L83").Value = 1341.75
$ws.Range("M83").Value = 4453.42855
$ws.Range("N83").Value = -11325.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 83336950
$ws.Range("I62").Value = 83336950
$ws.Range("K62").Value = 83336950
$ws.Range("M62").Value = -83336326
$ws.Range("H65").Value = 83336950
$ws.Range("I65").Value = 83336950
$ws.Range("K65").Value = 416684750
$ws.Range("M65").Value = -416681630
$ws.Range("H99").Value = 13338296
$ws.Range("I99").Value = 28574504
$ws.Range("K99").Value = 28574504
$ws.Range("M99").Value = -28573006
$ws.Range("H107").Value = 647.14813
$ws.Range("I107").Value = 498.875
$ws.Range("J107").Value = 1833.3334
$ws.Range("K107").Value = 498.875
$ws.Range("L107").Value = 1833.3334
$ws.Range("M107").Value = 1421.125
$ws.Range("N107").Value = -5673.3334
$ws.Range("H126").Value = 13338296
$ws.Range("I126").Value = 28574504
$ws.Range("K126").Value = 85723512
$ws.Range("M126").Value = -85721042

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2964.7856
$ws.Range("I122").Value = 1940.7
$ws.Range("K122").Value = 5822.1
$ws.Range("M122").Value = -3372.1
$ws.Range("H132").Value = 3827.4707
$ws.Range("I132").Value = 2318.889
$ws.Range("J132").Value = 5524.625
$ws.Range("K132").Value = 6956.667
$ws.Range("L132").Value = 16573.875
$ws.Range("M132").Value = -4426.667
$ws.Range("N132").Value = -21633.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2060.389
$ws.Range("I61").Value = 1791.9286
$ws.Range("K61").Value = 1791.9286
$ws.Range("M61").Value = -1589.9286
$ws.Range("H82").Value = 1206.5927
$ws.Range("I82").Value = 732.5263
$ws.Range("J82").Value = 2332.5
$ws.Range("K82").Value = 732.5263
$ws.Range("L82").Value = 2332.5
$ws.Range("M82").Value = -371.5263
$ws.Range("N82").Value = -3054.5
$ws.Range("H85").Value = 1206.5927
$ws.Range("I85").Value = 732.5263
$ws.Range("J85").Value = 2332.5
$ws.Range("K85").Value = 732.5263
$ws.Range("L85").Value = 2332.5
$ws.Range("M85").Value = 515.4737
$ws.Range("N85").Value = -4828.5
$ws.Range("H93").Value = 2600
$ws.Range("I93").Value = 1433.3334
$ws.Range("K93").Value = 1433.3334
$ws.Range("M93").Value = -185.3334
$ws.Range("H113").Value = 2060.389
$ws.Range("I113").Value = 1791.9286
$ws.Range("K113").Value = 1791.9286
$ws.Range("M113").Value = 378.0714
$ws.Range("H122").Value = 5377.857
$ws.Range("I122").Value = 3129
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 9387
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -6937
$ws.Range("N122").Value = -37900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 770.03845
$ws.Range("I107").Value = 511.05884
$ws.Range("J107").Value = 1259.2222
$ws.Range("K107").Value = 1533.17652
$ws.Range("L107").Value = 3777.6666
$ws.Range("M107").Value = 386.82348
$ws.Range("N107").Value = -7617.6666
